$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 274, shifting existing rows 274-288 down to 275-289
$ws.Rows.Item(274).Insert()

# Fill in the new row 274 with the new data record
$ws.Cells.Item(274, 1).Value = 8
$ws.Cells.Item(274, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(274, 3).Value = "Coquimbo"
$ws.Cells.Item(274, 4).Value = 44753
$ws.Cells.Item(274, 5).Value = 4
$ws.Cells.Item(274, 6).Value = 100112012
$ws.Cells.Item(274, 7).Value = "Espinaca"
$ws.Cells.Item(274, 8).Value = "Sin especificar"
$ws.Cells.Item(274, 9).Value = "Primera"
$ws.Cells.Item(274, 10).Value = 2600
$ws.Cells.Item(274, 11).Value = 500
$ws.Cells.Item(274, 12).Value = 600
$ws.Cells.Item(274, 13).Value = 550
$ws.Cells.Item(274, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(274, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(274, 16).Value = 1100
$ws.Cells.Item(274, 17).Value = 0.5
$ws.Cells.Item(274, 18).Value = "Hortaliza"

# Match the date style used by column D in the template rows (numFmtId 165)
$ws.Cells.Item(274, 4).NumberFormat = $ws.Cells.Item(275, 4).NumberFormat
